$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.489.78'
$ws.Range("E2").Value = '  -0.71%  '
$ws.Range("D3").Value = '1.627.26'
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.50%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.77'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.496'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.31%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.01'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.34%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.249'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.68%  '
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.98'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0839'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.80%  '
$ws.Range("D12").Value = '1.854.78'
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").Value = '1.610.12'
$ws.Range("E13").Value = '  -1.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.11'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.521'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.86'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.45%  '
$ws.Range("D17").Value = '26.500.79'
$ws.Range("E17").Value = '  -0.59%  '
$ws.Range("E18").Value = '  +0.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '214.91'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.01'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.30'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("E22").Value = '  +1.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.31'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.98'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +4.66%  '
$ws.Range("E25").Value = '  +2.11%  '
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("E27").Value = '  -0.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.83'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.75%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.54'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0507'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.57%  '
$ws.Range("E31").Value = '  -1.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.31'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.84%  '
$ws.Range("E33").Value = '  -0.34%  '
$ws.Range("E34").Value = '  -0.31%  '
$ws.Range("D35").Value = '1.220.67'
$ws.Range("E35").Value = '  +4.73%  '
$ws.Range("E36").Value = '  -1.50%  '
$ws.Range("E37").Value = '  +5.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.01'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.797'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.501'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("E41").Value = '  -2.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.795'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.33%  '
$ws.Range("D44").Value = '1.764.83'
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.83'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.66%  '
$ws.Range("E46").Value = '  +1.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.80'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.50%  '
$ws.Range("E48").Value = '  -1.55%  '
$ws.Range("E49").Value = '  -0.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.63'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.408'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.08%  '
